$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Header: new year column 2019 in AO1
$ws.Range("AO1").Value = 2019

# Fill AO2:AO218 with the ".." placeholder (missing-data marker) used
# throughout this sheet for years without reported data.
$ws.Range("AO2:AO218").Value = ".."

# Reflect the author's post-edit selection/scroll state.
$ws.Range("AO2:AO218").Select()
$ws.Application.ActiveWindow.ScrollColumn = 19
